$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the text-valued columns (Percentage, Target_Distribution) stay as text
# rather than being auto-converted by Excel (e.g. "59.8%" -> percentage number).
$ws.Range("C2:C4").NumberFormat = "@"
$ws.Range("E2:E4").NumberFormat = "@"

# Row 2 (Train)
$ws.Range("B2").Value = 250
$ws.Range("C2").Value = "59.8%"
$ws.Range("D2").Value = 252
$ws.Range("E2").Value = "no: 63.60%, yes: 36.40%"

# Row 3 (Val)
$ws.Range("B3").Value = 84
$ws.Range("C3").Value = "20.1%"
$ws.Range("D3").Value = 74
$ws.Range("E3").Value = "no: 64.29%, yes: 35.71%"

# Row 4 (Test)
$ws.Range("B4").Value = 84
$ws.Range("C4").Value = "20.1%"
$ws.Range("D4").Value = 88
$ws.Range("E4").Value = "no: 63.10%, yes: 36.90%"
